$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)
$s = $ftr.Range.InlineShapes.Item(1)
$rng = $s.Range
$ft = $rng.FormattedText
Write-Output ("FormattedText class: " + $ft.GetType())
$endRng = $d.Content
$endRng.Collapse(0)
$endRng.FormattedText = $ft
Write-Output "done"
